$d = $word.ActiveDocument

# The diff appends, after the existing "Week 18" entry, a new "Week 19"
# blog entry: two blank paragraphs, a Heading2 "Week 19" paragraph, one
# more blank paragraph, and the body paragraph for the new entry.

# Paragraph: blank
$d.Paragraphs.Last.Range.InsertParagraphAfter()

# Paragraph: blank
$d.Paragraphs.Last.Range.InsertParagraphAfter()

# Paragraph: "Week 19" heading. Set the text before applying the Heading2
# style so the style doesn't carry over into the paragraphs inserted
# after it.
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$headingPara = $d.Paragraphs.Last
$headingPara.Range.Text = "Week 19"

# Paragraph: blank
$d.Paragraphs.Last.Range.InsertParagraphAfter()

# Paragraph: the new week's body text
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "I reached out to Mr Judhi to update him on my project progress. I highlighted the progress on the AI model and the challenges of hosting the model on the Railway cloud server. I had no other issues to discuss with him, so we had a short encounter. Next week is reading week, so students and staff will not be in school. I will try to meet with him the week after."

# Apply the heading style now that every paragraph has been created.
$headingPara.Style = "Heading2"
